# Generate Report for Handback
#
# Populates the "Latest Target File" (F) / "Latest Handback File" (G) columns
# on the per-locale sheets, flips the Status from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears (Overview +
# per-locale sheets), and stamps fresh "Latest Handback DateTime" (H) values.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/da13284d45f18ef3a992b8727d432a18cc4e8b4e/e2e/"
$file1 = "11f943d4-e580-437c-8da3-f377f2e05c15"
$file2 = "bb898a2c-ec1c-4f15-987d-4e981b05dd91"

# ---------------------------------------------------------------------
# Overview sheet: just the status text needs to flip (columns B & C).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), ($mdBase + $file1 + ".md"), "", "", ($file1 + ".md"))
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c7cae3eadf7f18b3c548cd26d9555be48b7c6f7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/11f943d4-e580-437c-8da3-f377f2e05c15.3f056b18b7de5a52b2f2b3434c918479e351e156.zh-cn.xlf", "", "", "11f943d4-e580-437c-8da3-f377f2e05c15.3f056b18b7de5a52b2f2b3434c918479e351e156.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), ($mdBase + $file2 + ".md"), "", "", ($file2 + ".md"))
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c7cae3eadf7f18b3c548cd26d9555be48b7c6f7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bb898a2c-ec1c-4f15-987d-4e981b05dd91.7f7a0c4fed9b66d3029f48ed7d568ae715cbcbe9.zh-cn.xlf", "", "", "bb898a2c-ec1c-4f15-987d-4e981b05dd91.7f7a0c4fed9b66d3029f48ed7d568ae715cbcbe9.zh-cn.xlf")

$wsZh.Range("H2").Value = "2016-03-18 22:49:44"
$wsZh.Range("H3").Value = "2016-03-18 22:49:44"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), ($mdBase + $file1 + ".md"), "", "", ($file1 + ".md"))
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7adf40773f59f8bce68bcace7ba22d63e77bcecb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/11f943d4-e580-437c-8da3-f377f2e05c15.3f056b18b7de5a52b2f2b3434c918479e351e156.de-de.xlf", "", "", "11f943d4-e580-437c-8da3-f377f2e05c15.3f056b18b7de5a52b2f2b3434c918479e351e156.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), ($mdBase + $file2 + ".md"), "", "", ($file2 + ".md"))
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7adf40773f59f8bce68bcace7ba22d63e77bcecb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bb898a2c-ec1c-4f15-987d-4e981b05dd91.7f7a0c4fed9b66d3029f48ed7d568ae715cbcbe9.de-de.xlf", "", "", "bb898a2c-ec1c-4f15-987d-4e981b05dd91.7f7a0c4fed9b66d3029f48ed7d568ae715cbcbe9.de-de.xlf")

$wsDe.Range("H2").Value = "2016-03-18 22:49:49"
$wsDe.Range("H3").Value = "2016-03-18 22:49:49"
